$d = $word.ActiveDocument

# Locate the final paragraph (the one ending in "...every time he passes." that
# currently also carries the _GoBack bookmark).
$targetIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replace the whole paragraph's range with: the same two original runs (now
# without the _GoBack bookmark, which is being relocated below) followed by
# the newly authored paragraphs. A single InsertXML call lets us reproduce
# the exact paragraph/run/bookmark layout from the target revision.
$newXml = @"
<w:p $wns><w:r><w:t>Keeps restarting when he dies, in the start tried to help everyone, in the end he became corrupted by the repeated restarts, leading to logical leaps making him think everyone can be saved if he controls magic and the lay-lines.</w:t></w:r><w:r><w:t xml:space="preserve"> All this started when he was at death’s door by his already dead comrades, he wishes with his last breath to be able to save everyone as he closed his eyes for the final rest. He discovers that he has been sent back in time to the week prior when the group was in a bar, with no recollection of why or how. He discovered that this happens every time he dies, no one believes him when he tells them, but they stay together, the rest of the group accepting his recollection of the future as one of his quirks that should just be ignored. He keeps dying to save all his comrades not knowing the eldritch monstrosity lurking outside of time slowly corrupting him every time he passes.</w:t></w:r></w:p>
<w:p $wns></w:p>
<w:p $wns></w:p>
<w:p $wns><w:pPr><w:pStyle w:val="Overskrift1"/></w:pPr><w:r><w:t>Interesting characters:</w:t></w:r></w:p>
<w:p $wns><w:proofErr w:type="spellStart"/><w:r><w:t>Gaster</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-like character, spread throughout time, can noticed if people roll sufficiently high, is everywhere always… just hidden in the weave.</w:t></w:r><w:r><w:t xml:space="preserve"> Is incomprehensible, sounds like garbled mechanical speech of every language at the same time.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p $wns></w:p>
<w:p $wns><w:pPr><w:pStyle w:val="Overskrift1"/></w:pPr><w:r><w:t>Interesting items</w:t></w:r></w:p>
<w:p $wns><w:r><w:t>Items that let’s you look at the weave and adjust fate.</w:t></w:r></w:p>
"@

$r.InsertXML($newXml)
